{"js": "// The document's trailing paragraph is currently empty. Fill it with the\n// new question text, then insert a brand-new paragraph right after it\n// that holds the \"A1 :\" answer stub \u2014 matching the target diff.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst qParagraph = paragraphs.items[paragraphs.items.length - 1];\n\nconst qRange = qParagraph.insertText(\n  \"Q: How to calculate the accuracy of a binary classification algorithm using confusion matrix?\",\n  Word.InsertLocation.replace\n);\nqRange.font.name = \"Arial\";\nqRange.font.size = 14;\nqRange.font.color = \"#1D1C1D\";\nqRange.font.nameFarEast = \"Times New Roman\";\nqRange.font.nameBidirectional = \"Arial\";\nqRange.font.sizeBidirectional = 14;\nqRange.languageIdFarEast = \"tr-TR\";\n\nconst aParagraph = qParagraph.insertParagraph(\"A1 :\", Word.InsertLocation.after);\naParagraph.alignment = Word.Alignment.justified;\n\nconst aRange = aParagraph.getRange();\naRange.font.name = \"Arial\";\naRange.font.size = 14;\naRange.font.color = \"#1D1C1D\";\naRange.font.nameFarEast = \"Times New Roman\";\naRange.font.nameBidirectional = \"Arial\";\naRange.font.sizeBidirectional = 14;\naRange.languageIdFarEast = \"tr-TR\";\n\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# The last (trailing) paragraph in the document body is currently empty.\n# Fill it with the new question text, then append a brand-new paragraph\n# right after it containing the \"A1 :\" answer stub.\n\n$count = $d.Paragraphs.Count\n$qPara = $d.Paragraphs.Item($count)\n$qRange = $qPara.Range\n$qRange.Text = \"Q: How to calculate the accuracy of a binary classification algorithm using confusion matrix?\"\n$qRange.Font.NameAscii = \"Arial\"\n$qRange.Font.NameFarEast = \"Times New Roman\"\n$qRange.Font.NameBi = \"Arial\"\n$qRange.Font.Name = \"Arial\"\n$qRange.Font.Size = 14\n$qRange.Font.SizeBi = 14\n$qRange.Font.Color = 1907741\n$qRange.LanguageIDFarEast = \"tr-TR\"\n\n$qRange.InsertParagraphAfter()\n\n$newCount = $d.Paragraphs.Count\n$aPara = $d.Paragraphs.Item($newCount)\n$aRange = $aPara.Range\n$aRange.Text = \"A1 :\"\n$aRange.Font.NameAscii = \"Arial\"\n$aRange.Font.NameFarEast = \"Times New Roman\"\n$aRange.Font.NameBi = \"Arial\"\n$aRange.Font.Name = \"Arial\"\n$aRange.Font.Size = 14\n$aRange.Font.SizeBi = 14\n$aRange.Font.Color = 1907741\n$aRange.LanguageIDFarEast = \"tr-TR\"\n"}
